$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert the E2:E8 "squared" formulas into a single shared formula ---
# (Re-entering the same formula across the whole range makes Excel store it
# as one shared formula, anchored at E2, instead of 7 separate <f> entries.)
$ws.Range("E2:E8").Formula = "=A2*A2"

# --- New scenario rows (17-20) with a couple more N-body data points ---
# C17:C19 body counts; B20 typed directly (no formula for the last row's B).
$ws.Range("C17").Value = 3
$ws.Range("C18").Value = 4
$ws.Range("C19").Value = 5
$ws.Range("B20").Value = 28

# N Squared column (A) derived back from the total time column.
$ws.Range("A17").Formula = "=SQRT(E17)"
$ws.Range("A18:A20").Formula = "=SQRT(E18)"

# Accel. time column (B), scaled off C13 (time per frame) and body count.
$ws.Range("B17").Formula = "=C`$13*C17 - 15"
$ws.Range("B18:B19").Formula = "=C`$13*C18 - 15"

# Total time column (E), scaled proportionally from the row-7 baseline.
$ws.Range("E17").Formula = "=B17/B`$7*E`$7"
$ws.Range("E18:E20").Formula = "=B18/B`$7*E`$7"

# New number format (plain thousands separator, no decimals) for the new
# "Total Time" column cells.
$ws.Range("E17:E20").NumberFormat = "#,##0"

# --- Update the active selection to the newly added cells ---
$ws.Range("E19:E20").Select() | Out-Null
